$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text edits (preserve rich-text runs, edit substrings in place) ---
$ws.Range("A8").Characters(21, 2).Text = "25"
$ws.Range("C9").Characters(27, 9).Text = "6/19/2023"
$ws.Range("C9").Characters(47, 9).Text = "6/25/2023"

# --- Data table edits ---
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 4
$ws.Range("K15").Value = 25
$ws.Range("M15").Value = -37.5
$ws.Range("I14").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 2
$ws.Range("H16").Value = -71.428571428571
$ws.Range("I16").Value = 27
$ws.Range("J16").Value = 30
$ws.Range("K16").Value = -10
$ws.Range("L16").Value = 35
$ws.Range("M16").Value = -47.058823529411
$ws.Range("N16").Value = -80.985915492957
$ws.Range("C17").Value = 3
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 28.571428571428
$ws.Range("I17").Value = 89
$ws.Range("K17").Value = 111.904761904762
$ws.Range("L17").Value = 50.847457627118
$ws.Range("M17").Value = 32.835820895522
$ws.Range("N17").Value = -34.558823529411
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 16.666666666666
$ws.Range("I18").Value = 46
$ws.Range("J18").Value = 34
$ws.Range("K18").Value = 35.294117647058
$ws.Range("L18").Value = 21.052631578947
$ws.Range("M18").Value = -46.511627906976
$ws.Range("N18").Value = -92.933947772657
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -30
$ws.Range("F19").Value = 31
$ws.Range("G19").Value = 29
$ws.Range("H19").Value = 6.896551724137
$ws.Range("I19").Value = 229
$ws.Range("J19").Value = 148
$ws.Range("K19").Value = 54.729729729729
$ws.Range("L19").Value = 47.741935483871
$ws.Range("M19").Value = 22.459893048128
$ws.Range("N19").Value = -43.872549019607
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -66.666666666666
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = -46.153846153846
$ws.Range("I20").Value = 47
$ws.Range("J20").Value = 62
$ws.Range("K20").Value = -24.193548387096
$ws.Range("L20").Value = 104.347826086957
$ws.Range("M20").Value = 14.634146341463
$ws.Range("N20").Value = -96.463506395786
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = -6.25
$ws.Range("F21").Value = 58
$ws.Range("G21").Value = 64
$ws.Range("H21").Value = -9.375
$ws.Range("I21").Value = 444
$ws.Range("J21").Value = 321
$ws.Range("K21").Value = 38.317757009345
$ws.Range("L21").Value = 47.50830564784
$ws.Range("M21").Value = 0.90909090909
$ws.Range("N21").Value = -83.426651735722
$ws.Range("I14").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("D23").Value = 2
$ws.Range("K14").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("E23").Value = -100
$ws.Range("G23").Value = 4
$ws.Range("J23").Value = 12
$ws.Range("K23").Value = 66.666666666666
$ws.Range("M23").Value = 150
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = 111.111111111111
$ws.Range("F24").Value = 80
$ws.Range("G24").Value = 49
$ws.Range("H24").Value = 63.265306122449
$ws.Range("I24").Value = 541
$ws.Range("J24").Value = 337
$ws.Range("K24").Value = 60.53412462908
$ws.Range("L24").Value = 132.188841201717
$ws.Range("M24").Value = -32.961586121437
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 21
$ws.Range("G25").Value = 26
$ws.Range("H25").Value = -19.230769230769
$ws.Range("I25").Value = 148
$ws.Range("J25").Value = 168
$ws.Range("K25").Value = -11.904761904761
$ws.Range("L25").Value = 29.824561403508
$ws.Range("M25").Value = -47.703180212014
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 11
$ws.Range("K26").Value = -27.272727272727
$ws.Range("C27").Value = 1
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 25
$ws.Range("I27").Value = 19
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 35.714285714285
$ws.Range("F28").NumberFormat = "@"
$ws.Range("F28").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("F29").NumberFormat = "@"
$ws.Range("F29").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("I14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D30").Value = 1
$ws.Range("K14").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value = -100
$ws.Range("J30").Value = 5
$ws.Range("K30").Value = -20
